# Add a new "Save" column (H) to the s_vals sheet, matching the header
# style already used by the other header cells (B1:G1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H1: new header label, styled like the existing headers (bold, bordered,
# centered) by copying G1's format.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# H2: data value for the new column, plain numeric cell (no special style).
$ws.Range("H2").Value = 0
